$d = $word.ActiveDocument

# Helper: locate the paragraph containing $searchText (a short, unique
# substring of the paragraph we want to target) and insert a brand new
# paragraph right after it, filled with $newText and formatted italic.
function Add-ItalicParagraphAfterFind($doc, $searchText, $newText) {
    $rng = $doc.Content.Duplicate
    $found = $rng.Find.Execute($searchText, $true, $false, $false, $false, $false,
                                $true, 1, $false, "", 0)
    if (-not $found) {
        return
    }
    $para = $rng.Paragraphs(1)
    $paraRange = $para.Range
    $paraRange.InsertParagraphAfter() | Out-Null
    $insertStart = $paraRange.End
    $newRange = $doc.Range($insertStart, $insertStart)
    $newRange.InsertAfter($newText)
    $endPos = $insertStart + $newText.Length
    $textRange = $doc.Range($insertStart, $endPos)
    $textRange.Font.Italic = $true
}

# 1. Update activation date: 15/07/2016 -> 01/01/2023
$d.Content.Find.Execute(
    "Ativação: 15/07/2016", $true, $false, $false, $false, $false,
    $true, 1, $false, "Ativação: 01/01/2023", 2) | Out-Null

# 2. Insert English (italic) translation after the "Objetivos" body paragraph
Add-ItalicParagraphAfterFind $d `
    "Introdução às funções de variáveis complexas e suas aplicações." `
    "Introduction to complex variable functions and their applications. Present differential equations of interest in physical engineering and develop solution techniques, verifying properties and resolution methods. Study of special functions in Physical Engineering."

# 3. Insert English (italic) translation after the "Programa resumido" body paragraph
Add-ItalicParagraphAfterFind $d `
    "Funções de uma variável complexa. Função delta." `
    "Functions of a complex variable. Delta function. Partial differential equations in physical engineering: solution methods, solving boundary value problems, applications. Fourier Series and Integral Transforms. Special functions."

# 4. Insert English (italic) translation after the "Programa" body paragraph
Add-ItalicParagraphAfterFind $d `
    "Funções de uma variável complexa: séries infinitas" `
    "Functions of a complex variable: infinite series, analytical functions, Cauchy Riemann conditions, boundary integrals, Cauchy's theorem, residue theorem, Delta function. Laplace equation, diffusion equation (of heat), wave equation (vibrating string), Fourier series, Fourier and Laplace integral transforms. Special functions: Legendre Polynomials, Spherical Harmonics, Bessel Functions."

# 5. Update evaluation criterion formula: (P1 + 2P2)/3 -> (P1 + P2)/2
$d.Content.Find.Execute(
    "Duas provas escritas: conceitos P1 e P2. Conceito Final = (P1 + 2P2)/3",
    $true, $false, $false, $false, $false,
    $true, 1, $false,
    "Duas provas escritas: conceitos P1 e P2. Conceito Final = (P1 + P2)/2", 2) | Out-Null
